# Updated by script on 11-11-2023 20:45
# Re-applies the upstream re-scrape: several match rows had their
# home/away/odds/url payload (columns F:V) re-ordered within the same
# matchday, and three new matches were appended at the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bufferRow = 1000

function Swap-Rows($rowA, $rowB) {
    # Stash row A's F:V payload in the buffer row, move B's payload into A,
    # then move the stashed A payload into B. Using PasteSpecial(xlPasteAll)
    # (-4104) keeps values + formats + number formats identical to a plain
    # cell move.
    $ws.Range("F$rowA`:V$rowA").Copy() | Out-Null
    $ws.Range("F$bufferRow`:V$bufferRow").PasteSpecial(-4104) | Out-Null

    $ws.Range("F$rowB`:V$rowB").Copy() | Out-Null
    $ws.Range("F$rowA`:V$rowA").PasteSpecial(-4104) | Out-Null

    $ws.Range("F$bufferRow`:V$bufferRow").Copy() | Out-Null
    $ws.Range("F$rowB`:V$rowB").PasteSpecial(-4104) | Out-Null

    $ws.Range("F$bufferRow`:V$bufferRow").ClearContents()
}

# Simple pairwise swaps (each pair exchanges its match payload, columns F:V,
# while index/country/tournament/season/date in A:E stay put).
Swap-Rows 12 13
Swap-Rows 14 15
Swap-Rows 23 24
Swap-Rows 29 30
Swap-Rows 35 37
Swap-Rows 38 39
Swap-Rows 53 55
Swap-Rows 89 90
Swap-Rows 92 94
Swap-Rows 96 97
Swap-Rows 101 102

# 3-way rotation: new25 = old26, new26 = old28, new28 = old25
$ws.Range("F25:V25").Copy() | Out-Null
$ws.Range("F$bufferRow`:V$bufferRow").PasteSpecial(-4104) | Out-Null

$ws.Range("F26:V26").Copy() | Out-Null
$ws.Range("F25:V25").PasteSpecial(-4104) | Out-Null

$ws.Range("F28:V28").Copy() | Out-Null
$ws.Range("F26:V26").PasteSpecial(-4104) | Out-Null

$ws.Range("F$bufferRow`:V$bufferRow").Copy() | Out-Null
$ws.Range("F28:V28").PasteSpecial(-4104) | Out-Null

$ws.Range("F$bufferRow`:V$bufferRow").ClearContents()

# Append three brand-new matches (rows 109-111) at the bottom, copying the
# row-108 formatting down first so the style indices (bold/bordered index
# column, datetime-formatted date column) line up with the rest of the
# sheet.
$ws.Range("A108:V108").Copy() | Out-Null
$ws.Range("A109:V111").PasteSpecial(-4122) | Out-Null

$ws.Range("A109").Value = 108
$ws.Range("B109").Value = "spain"
$ws.Range("C109").Value = "primera-rfef-group-2"
$ws.Range("D109").Value = "2023-2024"
$ws.Range("E109").Value = 45241.66666666666
$ws.Range("F109").Value = "San Fernando"
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = "Antequera"
$ws.Range("I109").Value = 1
$ws.Range("J109").Value = 2.67
$ws.Range("K109").Value = "09/11/2023 09:13"
$ws.Range("L109").Value = 2.73
$ws.Range("M109").Value = "11/11/2023 15:54"
$ws.Range("N109").Value = 3.01
$ws.Range("O109").Value = "09/11/2023 09:13"
$ws.Range("P109").Value = 2.91
$ws.Range("Q109").Value = "11/11/2023 15:33"
$ws.Range("R109").Value = 2.59
$ws.Range("S109").Value = "09/11/2023 09:13"
$ws.Range("T109").Value = 2.88
$ws.Range("U109").Value = "11/11/2023 15:54"
$ws.Range("V109").Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/cd-san-fernando-antequera/O0hpGA7o/"

$ws.Range("A110").Value = 109
$ws.Range("B110").Value = "spain"
$ws.Range("C110").Value = "primera-rfef-group-2"
$ws.Range("D110").Value = "2023-2024"
$ws.Range("E110").Value = 45241.75
$ws.Range("F110").Value = "Atl. Madrid B"
$ws.Range("G110").Value = 1
$ws.Range("H110").Value = "Recreativo Huelva"
$ws.Range("I110").Value = 2
$ws.Range("J110").Value = 1.7
$ws.Range("K110").Value = "09/11/2023 09:13"
$ws.Range("L110").Value = 1.85
$ws.Range("M110").Value = "11/11/2023 16:32"
$ws.Range("N110").Value = 3.43
$ws.Range("O110").Value = "09/11/2023 09:13"
$ws.Range("P110").Value = 3.32
$ws.Range("Q110").Value = "11/11/2023 16:32"
$ws.Range("R110").Value = 4.72
$ws.Range("S110").Value = "09/11/2023 09:13"
$ws.Range("T110").Value = 4.65
$ws.Range("U110").Value = "11/11/2023 16:32"
$ws.Range("V110").Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/atl-madrid-recreativo-huelva/tMbRakN4/"

$ws.Range("A111").Value = 110
$ws.Range("B111").Value = "spain"
$ws.Range("C111").Value = "primera-rfef-group-2"
$ws.Range("D111").Value = "2023-2024"
$ws.Range("E111").Value = 45241.83333333334
$ws.Range("F111").Value = "Algeciras"
$ws.Range("G111").Value = 0
$ws.Range("H111").Value = "CF Intercity"
$ws.Range("I111").Value = 1
$ws.Range("J111").Value = 2.15
$ws.Range("K111").Value = "09/11/2023 09:13"
$ws.Range("L111").Value = 2.38
$ws.Range("M111").Value = "11/11/2023 19:59"
$ws.Range("N111").Value = 2.88
$ws.Range("O111").Value = "09/11/2023 09:13"
$ws.Range("P111").Value = 2.88
$ws.Range("Q111").Value = "11/11/2023 19:59"
$ws.Range("R111").Value = 3.45
$ws.Range("S111").Value = "09/11/2023 09:13"
$ws.Range("T111").Value = 3.44
$ws.Range("U111").Value = "11/11/2023 19:59"
$ws.Range("V111").Value = "https://www.betexplorer.com/football/spain/primera-rfef-group-2/algeciras-cf-intercity/fP2Vb9xB/"
